$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.99999999106231852
$ws.Range("A2").Value = 0.99380268995620791
$ws.Range("A3").Value = 0.97012793306358747
$ws.Range("A4").Value = 0.95855279300161833
$ws.Range("A5").Value = 0.94733476313692422
$ws.Range("A6").Value = 0.91996954567171085
$ws.Range("A7").Value = 0.9133658159502942
$ws.Range("A8").Value = 0.90422763641587467
$ws.Range("A9").Value = 0.89247162225527776
$ws.Range("A10").Value = 0.88163632778184919
$ws.Range("A11").Value = 0.88003940179238882
$ws.Range("A12").Value = 0.87725816368647469
$ws.Range("A13").Value = 0.86596628153058297
$ws.Range("A14").Value = 0.86179792791249188
$ws.Range("A15").Value = 0.85920595339442851
$ws.Range("A16").Value = 0.8566991068029437
$ws.Range("A17").Value = 0.85299085015053455
$ws.Range("A18").Value = 0.85188193320077521
$ws.Range("A19").Value = 0.99046379187690048
$ws.Range("A20").Value = 0.98334672854162208
$ws.Range("A21").Value = 0.981948247972807
$ws.Range("A22").Value = 0.98068374182005202
$ws.Range("A23").Value = 0.96309215486516919
$ws.Range("A24").Value = 0.95007026353431634
$ws.Range("A25").Value = 0.94361310543432286
$ws.Range("A26").Value = 0.91791781943904294
$ws.Range("A27").Value = 0.91530784471832882
$ws.Range("A28").Value = 0.9037796364800148
$ws.Range("A29").Value = 0.89591123450401211
$ws.Range("A30").Value = 0.89319900216956405
$ws.Range("A31").Value = 0.89487075084686774
$ws.Range("A32").Value = 0.89198287775795404
$ws.Range("A33").Value = 0.88989030637329136
